$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: Activity 2 paragraph - "... to collect up-cycl" -> "... to collect
# up-cycle and down-cycle data sets." (also drops the spell-check proofErr
# markers that wrapped "cycl" in the original).
# ---------------------------------------------------------------------------

function Get-ParaIndexByText($needle) {
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        if ($d.Paragraphs($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

$idx1 = Get-ParaIndexByText "to collect up-cycl"
$prev1 = $d.Paragraphs($idx1 - 1)

$text1 = "Depress the pedal through a range of known inputs while recording the corresponding output signal voltages. A suggested procedure for determining the known input angle is to measure the vertical position of a reference point on the moving pedal with respect to the fixed base. Trigonometry can be used to calculate the pedal angle. Record the measured distances, the calculated angle, and the output voltage for each position. Repeat the process while increasing the input and decreasing the input to collect up-cycle and down-cycle data sets."

$prev1.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs($idx1)
$newPara1.Range.Text = $text1

$oldPara1 = $d.Paragraphs($idx1 + 1)
$oldRange1 = $d.Range($oldPara1.Range.Start, $oldPara1.Range.End)
$oldRange1.Delete()

# ---------------------------------------------------------------------------
# Edit 2: Activity 3 paragraph - move the "Show both sets..." sentence so it
# comes right before "indicate which points come from increasing...", drop
# the trailing period after "background document", and relocate the
# _GoBack bookmark from the end of this paragraph to the end of the second
# blank paragraph that follows it (right before "Activity 4").
# ---------------------------------------------------------------------------

$idx2 = Get-ParaIndexByText "Plot the collected calibration data"
$prev2 = $d.Paragraphs($idx2 - 1)

$text2 = "Plot the collected calibration data with the known reference data on the x axis and the measured voltages from the sensor signal on the y axis. Show both sets of calibration data on the same figure, and indicate which points come from increasing input (up-cycle) and decreasing input (down-cycle).  Use linear least squares regression to calculate the static sensitivity and zero offset of the resulting calibration curve. The required equations are given in the background document "

$prev2.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs($idx2)
$newPara2.Range.Text = $text2

$oldPara2 = $d.Paragraphs($idx2 + 1)
$oldRange2 = $d.Range($oldPara2.Range.Start, $oldPara2.Range.End)
$oldRange2.Delete()

# The two blank BodyText paragraphs right after $idx2 are untouched by the
# deletion above; the bookmark now belongs at the end of the second one.
$blank2 = $d.Paragraphs($idx2 + 2)
$d.Bookmarks.Add("_GoBack", $blank2.Range)

Write-Output "done"
